$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $text into cell $addr while FORCING Excel to treat it as a literal
# text string (matches the workbook author authoring tool, which always stores these
# cells as inline/shared strings) rather than letting Excel auto-coerce number-looking
# text (e.g. "314.61", "1.010") into a numeric cell. We do this by writing the text as
# a formula result into a scratch cell, then Copy + PasteSpecial(values) it onto the
# destination - PasteSpecial of a computed string preserves its Text cell type.
$helper = $ws.Range("ZZ1")
function Set-TextValue([string]$targetAddr, [string]$text) {
    $escaped = $text -replace '"', '""'
    $helper.Formula = '="' + $escaped + '"'
    $helper.Copy()
    $ws.Range($targetAddr).PasteSpecial(-4163)
}

# Row 2
Set-TextValue "D2" '27.349.42'
Set-TextValue "E2" '  +1.58%  '

# Row 3
Set-TextValue "D3" '1.832.66'
Set-TextValue "E3" '  +1.07%  '

# Row 4
Set-TextValue "E4" '  +0.93%  '

# Row 5
Set-TextValue "D5" '314.61'
Set-TextValue "E5" '  +1.77%  '

# Row 6
Set-TextValue "E6" '  +0.80%  '

# Row 7
Set-TextValue "D7" '0.4750'
Set-TextValue "E7" '  +2.07%  '

# Row 8
Set-TextValue "D8" '0.3689'
Set-TextValue "E8" '  +0.95%  '

# Row 9
Set-TextValue "D9" '0.07461'

# Row 10
Set-TextValue "D10" '0.8858'
Set-TextValue "E10" '  +1.96%  '

# Row 11
Set-TextValue "D11" '20.43'
Set-TextValue "E11" '  +1.06%  '

# Row 12
Set-TextValue "D12" '1.892.02'
Set-TextValue "E12" '  +3.29%  '

# Row 13
Set-TextValue "D13" '0.07338'
Set-TextValue "E13" '  +3.42%  '

# Row 14
Set-TextValue "D14" '5.444'
Set-TextValue "E14" '  +1.48%  '

# Row 15
Set-TextValue "D15" '93.14'
Set-TextValue "E15" '  +2.17%  '

# Row 16
Set-TextValue "D16" '6.587'
Set-TextValue "E16" '  +1.34%  '

# Row 17
Set-TextValue "D17" '1.010'
Set-TextValue "E17" '  +0.67%  '

# Row 18
Set-TextValue "D18" '0.000008801'
Set-TextValue "E18" '  +1.14%  '

# Row 19
Set-TextValue "E19" '  +0.85%  '

# Row 20
Set-TextValue "D20" '27.615.22'
Set-TextValue "E20" '  +2.50%  '

# Row 21
Set-TextValue "D21" '14.78'
Set-TextValue "E21" '  +1.08%  '

# Row 22
Set-TextValue "D22" '5.303'
Set-TextValue "E22" '  +0.29%  '

# Row 24
Set-TextValue "D24" '2.106.26'
Set-TextValue "E24" '  +2.87%  '

# Row 25
Set-TextValue "D25" '1.908'
Set-TextValue "E25" '  +0.72%  '

# Row 26
Set-TextValue "D26" '151.97'
Set-TextValue "E26" '  +0.78%  '

# Row 27
Set-TextValue "D27" '18.63'
Set-TextValue "E27" '  +1.73%  '

# Row 28
Set-TextValue "D28" '2.143'
Set-TextValue "E28" '  +0.93%  '

# Row 29
Set-TextValue "D29" '5.240'
Set-TextValue "E29" '  -0.28%  '

# Row 30
Set-TextValue "D30" '117.33'
Set-TextValue "E30" '  +1.56%  '

# Row 31
Set-TextValue "D31" '0.09000'
Set-TextValue "E31" '  +1.30%  '

# Row 32
Set-TextValue "D32" '0.7574'
Set-TextValue "E32" '  +0.51%  '

# Row 33
Set-TextValue "D33" '1.177'
Set-TextValue "E33" '  +1.38%  '

# Row 34
Set-TextValue "D34" '4.545'
Set-TextValue "E34" '  +1.47%  '

# Row 35
Set-TextValue "D35" '2.950'
Set-TextValue "E35" '  +1.62%  '

# Row 36
Set-TextValue "D36" '1.011'
Set-TextValue "E36" '  +0.88%  '

# Row 37
Set-TextValue "D37" '1.102'
Set-TextValue "E37" '  +1.60%  '

# Row 38
Set-TextValue "D38" '0.05350'
Set-TextValue "E38" '  +1.25%  '

# Row 39
Set-TextValue "D39" '0.01954'
Set-TextValue "E39" '  +0.46%  '

# Row 40
Set-TextValue "D40" '2.980'
Set-TextValue "E40" '  +0.34%  '

# Row 41
Set-TextValue "E41" '  +0.80%  '

# Row 42
Set-TextValue "D42" '2.401'
Set-TextValue "E42" '  +4.39%  '

# Row 43
Set-TextValue "D43" '0.5318'
Set-TextValue "E43" '  +0.21%  '

# Row 44
Set-TextValue "D44" '0.1660'
Set-TextValue "E44" '  +0.43%  '

# Row 45
Set-TextValue "D45" '8.478'
Set-TextValue "E45" '  +0.70%  '

# Row 46
Set-TextValue "D46" '0.4913'
Set-TextValue "E46" '  +1.05%  '

# Row 47
Set-TextValue "D47" '10.53'
Set-TextValue "E47" '  +1.46%  '

# Row 48
Set-TextValue "B48" 'PaxDollar'
Set-TextValue "C48" 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D48" '1.011'
Set-TextValue "E48" '  +0.93%  '

# Row 49
Set-TextValue "B49" 'Quant'
Set-TextValue "C49" 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D49" '104.90'
Set-TextValue "E49" '  +1.90%  '

# Row 51
Set-TextValue "D51" '0.06301'
Set-TextValue "E51" '  +0.18%  '

# Clean up scratch cell so it does not linger in the saved workbook
$helper.Clear()
